$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the daily deaths value for the last existing row (row 74, column J)
$ws.Range("J74").Value = 1

# Copy the formatting of row 73 (which already has the "no border, odd-row"
# style) down into the new row 75, so the new cells get the right styles
# before the values are written (avoids general/text formatting issues).
$ws.Range("A73:J73").Copy()
$ws.Range("A75:J75").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Append a new data row (row 75) with the day's figures
$ws.Range("A75").Value = 43975
$ws.Range("B75").Value = 75016
$ws.Range("C75").Value = 256
$ws.Range("D75").Value = 1469
$ws.Range("E75").Value = 1
$ws.Range("F75").Value = 16
$ws.Range("G75").Value = 3
$ws.Range("H75").Value = 2
$ws.Range("I75").Value = 107
$ws.Range("J75").Value = 0

$ws.Range("J74").Select()
